# Auto-generated edit script: updates market-data derived value cells
# across the Hyperion_Profits workbook (columns H-N) per the scheduled runner diff.

$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2651108.8
$ws.Range("I76").Value = 4278446
$ws.Range("K76").Value = 4278446
$ws.Range("M76").Value = -4278131
$ws.Range("H79").Value = 2651108.8
$ws.Range("I79").Value = 4278446
$ws.Range("K79").Value = 4278446
$ws.Range("M79").Value = -4277354
$ws.Range("H98").Value = 1359.5358
$ws.Range("I98").Value = 1317.2963
$ws.Range("K98").Value = 1317.2963
$ws.Range("M98").Value = 180.7037
$ws.Range("H122").Value = 1359.5358
$ws.Range("I122").Value = 1317.2963
$ws.Range("K122").Value = 3951.8889
$ws.Range("M122").Value = -1501.8889
$ws.Range("H137").Value = 102590.836
$ws.Range("I137").Value = 163624.36
$ws.Range("K137").Value = 490873.08
$ws.Range("M137").Value = -488323.08
$ws.Range("H138").Value = 8685.241
$ws.Range("J138").Value = 9095.23
$ws.Range("L138").Value = 27285.69
$ws.Range("N138").Value = -37565.69

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16730.229
$ws.Range("I32").Value = 13530.712
$ws.Range("K32").Value = 13530.712
$ws.Range("M32").Value = -13243.712
$ws.Range("H45").Value = 5497024
$ws.Range("I45").Value = 7144032.5
$ws.Range("J45").Value = 6995.1665
$ws.Range("K45").Value = 7144032.5
$ws.Range("L45").Value = 6995.1665
$ws.Range("M45").Value = -7143655.5
$ws.Range("N45").Value = -7749.1665
$ws.Range("H61").Value = 6352.72
$ws.Range("I61").Value = 6517.0527
$ws.Range("J61").Value = 5832.3335
$ws.Range("K61").Value = 6517.0527
$ws.Range("L61").Value = 5832.3335
$ws.Range("M61").Value = -6305.0527
$ws.Range("N61").Value = -6256.3335
$ws.Range("H110").Value = 5578900
$ws.Range("I110").Value = 6945847
$ws.Range("J110").Value = 111113
$ws.Range("K110").Value = 6945847
$ws.Range("L110").Value = 111113
$ws.Range("M110").Value = -6943802
$ws.Range("N110").Value = -115203
$ws.Range("H122").Value = 9859083
$ws.Range("I122").Value = 14467241
$ws.Range("J122").Value = 1899535.9
$ws.Range("K122").Value = 43401723
$ws.Range("L122").Value = 5698607.699999999
$ws.Range("M122").Value = -43399273
$ws.Range("N122").Value = -5703507.699999999
$ws.Range("H136").Value = 6352.72
$ws.Range("I136").Value = 6517.0527
$ws.Range("J136").Value = 5832.3335
$ws.Range("K136").Value = 19551.1581
$ws.Range("L136").Value = 17497.0005
$ws.Range("M136").Value = -17001.1581
$ws.Range("N136").Value = -22597.0005

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 37306.25
$ws.Range("I26").Value = 28275.857
$ws.Range("K26").Value = 28275.857
$ws.Range("M26").Value = -27983.857
$ws.Range("H105").Value = 8938058
$ws.Range("I105").Value = 12502879
$ws.Range("J105").Value = 26004.5
$ws.Range("K105").Value = 12502879
$ws.Range("L105").Value = 26004.5
$ws.Range("M105").Value = -12501132
$ws.Range("N105").Value = -29498.5
$ws.Range("H107").Value = 11908686
$ws.Range("I107").Value = 11908686
$ws.Range("K107").Value = 11908686
$ws.Range("M107").Value = -11906766

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18646.95
$ws.Range("I31").Value = 2384.16
$ws.Range("J31").Value = 30967.242
$ws.Range("K31").Value = 2384.16
$ws.Range("L31").Value = 30967.242
$ws.Range("M31").Value = -2089.16
$ws.Range("N31").Value = -31557.242
$ws.Range("H34").Value = 18646.95
$ws.Range("I34").Value = 2384.16
$ws.Range("J34").Value = 30967.242
$ws.Range("K34").Value = 2384.16
$ws.Range("L34").Value = 30967.242
$ws.Range("M34").Value = -2182.16
$ws.Range("N34").Value = -31371.242

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10959838
$ws.Range("I4").Value = 11368990
$ws.Range("K4").Value = 34106970
$ws.Range("M4").Value = -34106858
$ws.Range("H40").Value = 27.35
$ws.Range("I40").Value = 27.2
$ws.Range("K40").Value = 108.8
$ws.Range("M40").Value = -39.8

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16680067
$ws.Range("I70").Value = 20004440
$ws.Range("K70").Value = 20004440
$ws.Range("M70").Value = -20004170
$ws.Range("H73").Value = 16680067
$ws.Range("I73").Value = 20004440
$ws.Range("K73").Value = 20004440
$ws.Range("M73").Value = -20003504
$ws.Range("H102").Value = 4654303.5
$ws.Range("I102").Value = 6173661.5
$ws.Range("K102").Value = 6173661.5
$ws.Range("M102").Value = -6172039.5
$ws.Range("H126").Value = 5154724
$ws.Range("I126").Value = 3032931.5
$ws.Range("J126").Value = 8337412.5
$ws.Range("K126").Value = 9098794.5
$ws.Range("L126").Value = 25012237.5
$ws.Range("M126").Value = -9096324.5
$ws.Range("N126").Value = -25017177.5
$ws.Range("H132").Value = 4706.5713
$ws.Range("I132").Value = 4706.5713
$ws.Range("K132").Value = 14119.7139
$ws.Range("M132").Value = -11589.7139

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 64476.715
$ws.Range("I22").Value = 81598.55
$ws.Range("J22").Value = 1696.6666
$ws.Range("K22").Value = 81598.55
$ws.Range("L22").Value = 1696.6666
$ws.Range("M22").Value = -81303.55
$ws.Range("N22").Value = -2286.6666
$ws.Range("H27").Value = 64476.715
$ws.Range("I27").Value = 81598.55
$ws.Range("J27").Value = 1696.6666
$ws.Range("K27").Value = 81598.55
$ws.Range("L27").Value = 1696.6666
$ws.Range("M27").Value = -81491.55
$ws.Range("N27").Value = -1910.6666
$ws.Range("H61").Value = 7408455
$ws.Range("I61").Value = 9260101
$ws.Range("J61").Value = 1871.6666
$ws.Range("K61").Value = 9260101
$ws.Range("L61").Value = 1871.6666
$ws.Range("M61").Value = -9259899
$ws.Range("N61").Value = -2275.6666
$ws.Range("H113").Value = 7408455
$ws.Range("I113").Value = 9260101
$ws.Range("J113").Value = 1871.6666
$ws.Range("K113").Value = 9260101
$ws.Range("L113").Value = 1871.6666
$ws.Range("M113").Value = -9257931
$ws.Range("N113").Value = -6211.6666

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6948401.5
$ws.Range("I81").Value = 7578892.5
$ws.Range("K81").Value = 15157785
$ws.Range("M81").Value = -15156724
$ws.Range("H84").Value = 6948401.5
$ws.Range("I84").Value = 7578892.5
$ws.Range("K84").Value = 75788925
$ws.Range("M84").Value = -75783621
$ws.Range("H133").Value = 78342.44500000001
$ws.Range("J133").Value = 78342.44500000001
$ws.Range("L133").Value = 78342.44500000001
$ws.Range("N133").Value = -88462.44500000001
